# Weekly fruit/vegetable price update:
# Insert one new "Sandia" (watermelon) price record at row 186 of the
# "Hortaliza, Vega Central Mapocho de Santiago - Sandia" sheet, pushing
# the existing rows 186-210 down to 187-211.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift existing data down by inserting a new blank row at 186
$ws.Rows.Item(186).Insert()

# Populate the new row with the new weekly record
$ws.Range("A186").Value = 9
$ws.Range("B186").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C186").Value = "Metropolitana"
$ws.Range("D186").Value = 44491
$ws.Range("E186").Value = 13
$ws.Range("F186").Value = 100112028
$ws.Range("G186").Value = "Sandia"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 250
$ws.Range("K186").Value = 800
$ws.Range("L186").Value = 900
$ws.Range("M186").Value = 850
$ws.Range("N186").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O186").Value = "Perú"
$ws.Range("P186").Value = 850
$ws.Range("Q186").Value = 1
$ws.Range("R186").Value = "Hortaliza"
